# Updates cryptos list price/volume(1h) cells (and the three rank swaps
# around Monero / PolygonEcosystemToken / EthereumClassic and
# Mantle / InjectiveProtocol) to match the refreshed GitHub Actions data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.778.48"
$ws.Range("E2").Value = "  -0.67%  "
$ws.Range("D3").Value = "2.464.37"
$ws.Range("E3").Value = "  -0.60%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "571.90"
$ws.Range("E5").Value = "  -0.81%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.37"
$ws.Range("E6").Value = "  +0.60%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("E8").Value = "  -1.54%  "
$ws.Range("E9").Value = "  -0.01%  "
$ws.Range("E10").Value = "  -0.15%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.18"
$ws.Range("E11").Value = "  -1.27%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.348"
$ws.Range("E12").Value = "  -1.60%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "29.04"
$ws.Range("E13").Value = "  +1.66%  "
$ws.Range("E14").Value = "  -1.77%  "
$ws.Range("D15").Value = "2.908.78"
$ws.Range("E15").Value = "  -0.72%  "
$ws.Range("D16").Value = "62.668.04"
$ws.Range("E16").Value = "  -0.82%  "
$ws.Range("D17").Value = "2.466.83"
$ws.Range("E17").Value = "  -0.51%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.67"
$ws.Range("E18").Value = "  -5.97%  "
$ws.Range("E19").Value = "  -2.60%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.35"
$ws.Range("E20").Value = "  +4.94%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.15"
$ws.Range("E21").Value = "  +0.45%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "321.59"
$ws.Range("E22").Value = "  -2.61%  "
$ws.Range("E23").Value = "  +0.01%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.23"
$ws.Range("E24").Value = "  +3.93%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "64.76"
$ws.Range("E25").Value = "  -2.17%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "641.58"
$ws.Range("E26").Value = "  -2.53%  "
$ws.Range("D27").Value = "2.592.12"
$ws.Range("E27").Value = "  -0.66%  "
$ws.Range("D28").Value = "0.0₃0963"
$ws.Range("E28").Value = "  -2.84%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.19%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.42"
$ws.Range("E30").Value = "  -4.26%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.90"
$ws.Range("E31").Value = "  -2.17%  "
$ws.Range("E32").Value = "  -2.14%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.133"
$ws.Range("E33").Value = "  +0.00%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.998"
$ws.Range("E34").Value = "  -0.06%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.50"
$ws.Range("E35").Value = "  -3.18%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.69"
$ws.Range("E36").Value = "  -2.12%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.38"
$ws.Range("E37").Value = "  -1.45%  "
$ws.Range("B38").Value = "PolygonEcosystemToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.365"
$ws.Range("E38").Value = "  -1.71%  "
$ws.Range("B39").Value = "EthereumClassic"
$ws.Range("C39").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.56"
$ws.Range("E39").Value = "  -1.15%  "
$ws.Range("B40").Value = "Monero"
$ws.Range("C40").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "149.91"
$ws.Range("E40").Value = "  -0.25%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.63"
$ws.Range("E41").Value = "  -1.90%  "
$ws.Range("E42").Value = "  -1.78%  "
$ws.Range("E43").Value = "  +0.03%  "
$ws.Range("D44").Value = "0.0₆0307"
$ws.Range("E44").Value = "  -3.97%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "154.04"
$ws.Range("E45").Value = "  -0.70%  "
$ws.Range("E46").Value = "  +1.03%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.56"
$ws.Range("E47").Value = "  -1.21%  "
$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "20.25"
$ws.Range("E48").Value = "  -0.85%  "
$ws.Range("B49").Value = "Mantle"
$ws.Range("C49").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.606"
$ws.Range("E49").Value = "  -0.50%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0509"
$ws.Range("E50").Value = "  -1.21%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0902"
$ws.Range("E51").Value = "  -1.70%  "
